# CIS Update April 1 - May 1
#
# Adds a new "April 2023" tab (the newest monthly CIS flag sheet) in front
# of "March 2023", built from a copy of the "March 2023" template so it
# keeps the same header/column formatting. The old 5-row flag list is
# trimmed down to the 3 rows relevant to April and repopulated with the
# April good/bad channel names.

$wb = $excel.ActiveWorkbook

# Work from the previous month's sheet (still the template for the new one).
$template = $wb.Worksheets.Item("March 2023")
[void]$template.Activate()
[void]$template.Range("A4").Select()

# Duplicate it immediately before itself -> new sheet lands in slot 1,
# "March 2023" (and everything after it) shifts right by one.
$template.Copy($template)

$newSheet = $wb.Worksheets.Item(1)
$newSheet.Name = "April 2023"

# Drop the extra "Mark as Bad" rows carried over from March (7 -> 4 rows).
$newSheet.Rows.Item(7).Delete()
$newSheet.Rows.Item(6).Delete()
$newSheet.Rows.Item(5).Delete()

# Fill in April's updates.
$newSheet.Range("A2").Value = "EBA_m05_c10_lowgain"
$newSheet.Range("A3").Value = "EBA_m09_c16_lowgain"
$newSheet.Range("B2").Value = "EBA_m13_c04_lowgain"
$newSheet.Range("B3").ClearContents()
$newSheet.Range("B4").ClearContents()

[void]$newSheet.Range("B7").Select()
